$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Katha no 461 100 7, Comprised of Converted Survey No 100 4,"
$ws.Range("A4").Value = "SY  No  99 1, Aduru Village,"
